$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Done. 1: ('Adam', 'Lind', 'Toronto Blue Jays'), 2: ('Yasiel', 'Puig', 'Los Angeles Dodgers')"
$ws.Range("E3").Value = "Done. 1: ('Robinson', 'Cano', 'Seattle Mariners'), 2: ('Adam', 'Lind', 'Toronto Blue Jays')"
